# Everest (intelligent traffic steering) deck — analytics-docs bug fix.
#
# Adds two new "Everest-sink (job)" round-rect shapes to slide 8, styled
# identically to the existing "Registry (Service)" shape (same gradient
# fill / outline / shadow / Mac "wrapping textbox" extension), just with a
# different corner-radius adjustment, position/size and label text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

function Get-ShapeById {
    param($slide, [int]$id)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

function Add-EverestSinkShape {
    param($slide, $template, [double]$left, [double]$top, [double]$width, [double]$height)

    $shp = $template.Duplicate()
    $shp.Name = "Everest-sink…"
    $shp.Left = $left
    $shp.Top = $top
    $shp.Width = $width
    $shp.Height = $height
    $shp.Adjustments.Item(1) = 0.39419

    # Setting text containing "(job)" directly via TextRange.Text splits
    # the paragraph into 3 runs around the parentheses. Seed a plain
    # placeholder (no punctuation) first, then overwrite that single run's
    # Text so the paragraph keeps exactly one <a:r>, matching the source.
    $tr = $shp.TextFrame.TextRange
    $tr.Text = "Everest-sink`rPLACEHOLDER"
    $tr.Paragraphs(2, 1).Runs(1, 1).Text = "(job)"

    return $shp
}

# Existing "Registry…" shape (id 701) is the closest visual template:
# same roundRect geometry, gradient fill, line, shadow and the Mac
# ma14:wrappingTextBoxFlag extension.
$template = Get-ShapeById $s 701

Add-EverestSinkShape $s $template 773.5829921259842 223.98551181102363 91.84984251968504 33.20952805905512 | Out-Null
Add-EverestSinkShape $s $template 776.3137795275591 301.45417792834644 91.84984251968504 33.20952805905512 | Out-Null
